$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "OrdinalEncoder rather than OneHotEncoder"
$ws.Range("D5").Value = "add back in stratification and use OrdinalEncoder rather than OneHotEncoder"
$ws.Range("C7").Value = "Adding logarithmic transformation to numeric features to help normalize data"

$ws.Range("F4").Select()
